# Applies the GitHub Actions "Updated cryptos list" price refresh to the
# cryptos worksheet: updates Price (D) and Volume(1h) (E) figures for most
# rows, and swaps the WrappedBTC/WrappedEther rows (17/18), including their
# Coin name, Link, Price and Volume values.
#
# Price (D) values that look like plain numbers ("591.07", "14.20", ...)
# are written through a brief "@" (Text) number-format so Excel keeps them
# as literal text (preserving trailing zeros / thousands-dot formatting)
# instead of silently converting them to floating point numbers; the cell
# style is then reset back to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.033.38"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "3.121.84"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.39%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.113.88"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "3.634.10"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.136.04"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.040.64"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.92%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -8.09%  "
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("D38").Value = "0.0₃0715"
$ws.Range("E38").Value = "  -4.62%  "
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "421.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.64%  "
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("E43").Value = "  -11.02%  "
$ws.Range("D44").Value = "2.887.38"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.266"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("E49").Value = "  -6.72%  "
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.67%  "
